$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the style of an existing header cell (AB1) onto the new headers
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122) # xlPasteFormats

# Fill in the team record data for each data row (2-48)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 29).Value = 95  # AC = Wins
    $ws.Cells.Item($r, 30).Value = 67  # AD = Losses
    $ws.Cells.Item($r, 31).Value = 0   # AE = Ties
}
